{"js": "// Update the worksheet date heading and the 25 division problems in the\n// practice table. Cells are addressed by their (row, column) position in\n// the table grid rather than by matching old text, since some target\n// values coincide with other (unrelated) source values elsewhere in the\n// table (e.g. \"23\u00f76=\" is both a pre-edit value and a post-edit value).\n\n// 1. Update the date/weekday heading paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.load(\"text\");\nawait context.sync();\n\nif (heading.text === \"2024-01-04 Thursday\") {\n  heading.insertText(\"2024-01-05 Friday\", Word.InsertLocation.replace);\n}\n\n// 2. Update the division problems table (5 populated rows of 5 cells,\n// interleaved with blank rows left for student work).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newRows = {\n  0: [\"11\u00f77=\", \"26\u00f76=\", \"28\u00f76=\", \"53\u00f75=\", \"43\u00f73=\"],\n  4: [\"72\u00f73=\", \"17\u00f75=\", \"83\u00f74=\", \"24\u00f78=\", \"33\u00f78=\"],\n  8: [\"48\u00f78=\", \"22\u00f74=\", \"45\u00f74=\", \"98\u00f76=\", \"56\u00f79=\"],\n  12: [\"22\u00f78=\", \"45\u00f76=\", \"28\u00f79=\", \"24\u00f75=\", \"59\u00f73=\"],\n  16: [\"41\u00f79=\", \"64\u00f76=\", \"68\u00f78=\", \"76\u00f79=\", \"23\u00f76=\"],\n};\n\nconst updatedValues = table.values.map((row, rowIndex) => {\n  const replacement = newRows[rowIndex];\n  return replacement ? replacement.slice() : row;\n});\n\ntable.values = updatedValues;\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and the 25 division problems in the\n# practice table. Cells are addressed by their (row, column) position in\n# the table grid rather than by matching old text, since some target\n# values coincide with other (unrelated) source values elsewhere in the\n# table (e.g. \"23\u00f76=\" is both a pre-edit value and a post-edit value);\n# blind document-wide Find/Replace could therefore mutate the wrong cell.\n\n$d = $word.ActiveDocument\n\n# 1. Update the date/weekday heading paragraph (first paragraph).\n# Note: Range.Text includes the trailing paragraph-mark character, so\n# trim it before comparing against the plain target string.\n$heading = $d.Paragraphs.Item(1).Range\n$headingText = $heading.Text.TrimEnd([char]13, [char]7)\nif ($headingText -eq \"2024-01-04 Thursday\") {\n    $heading.Text = \"2024-01-05 Friday\"\n}\n\n# 2. Update the division problems table (rows 1, 5, 9, 13, 17 contain the\n# 5-cell rows of problems; the other rows are blank work rows).\n$t = $d.Tables.Item(1)\n\n$newRows = @{\n    1  = @(\"11\u00f77=\", \"26\u00f76=\", \"28\u00f76=\", \"53\u00f75=\", \"43\u00f73=\")\n    5  = @(\"72\u00f73=\", \"17\u00f75=\", \"83\u00f74=\", \"24\u00f78=\", \"33\u00f78=\")\n    9  = @(\"48\u00f78=\", \"22\u00f74=\", \"45\u00f74=\", \"98\u00f76=\", \"56\u00f79=\")\n    13 = @(\"22\u00f78=\", \"45\u00f76=\", \"28\u00f79=\", \"24\u00f75=\", \"59\u00f73=\")\n    17 = @(\"41\u00f79=\", \"64\u00f76=\", \"68\u00f78=\", \"76\u00f79=\", \"23\u00f76=\")\n}\n\nforeach ($rowIndex in $newRows.Keys) {\n    $rowValues = $newRows[$rowIndex]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
